$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2325581395348837
$ws.Range("C2").Value = 0.4922480620155039
$ws.Range("J2").Value = 0.007751937984496124
$ws.Range("P2").Value = 0.186046511627907
$ws.Range("S2").Value = 0.08139534883720931
$ws.Range("B3").Value = 0.01526717557251908
$ws.Range("C3").Value = 0.01526717557251908
$ws.Range("J3").Value = 0.02290076335877863
$ws.Range("P3").Value = 0.7175572519083969
$ws.Range("S3").Value = 0.2290076335877863
$ws.Range("J4").Value = 0.125
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.225
$ws.Range("B6").Value = 0.03864734299516908
$ws.Range("D6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.05797101449275362
$ws.Range("J6").Value = 0.2801932367149759
$ws.Range("O6").Value = 0.03864734299516908
$ws.Range("Q6").Value = 0.1545893719806763
$ws.Range("R6").Value = 0.05797101449275362
$ws.Range("S6").Value = 0.3671497584541063
$ws.Range("B7").Value = 0.125
$ws.Range("D7").Value = 0.035
$ws.Range("F7").Value = 0.04
$ws.Range("J7").Value = 0.19
$ws.Range("O7").Value = 0.025
$ws.Range("Q7").Value = 0.13
$ws.Range("R7").Value = 0.06
$ws.Range("S7").Value = 0.395
$ws.Range("B8").Value = 0.0870712401055409
$ws.Range("D8").Value = 0.01846965699208443
$ws.Range("F8").Value = 0.06068601583113457
$ws.Range("J8").Value = 0.09762532981530343
$ws.Range("O8").Value = 0.0316622691292876
$ws.Range("Q8").Value = 0.1530343007915567
$ws.Range("R8").Value = 0.08970976253298153
$ws.Range("S8").Value = 0.4617414248021108
$ws.Range("B9").Value = 0.04938271604938271
$ws.Range("D9").Value = 0.01234567901234568
$ws.Range("F9").Value = 0.03703703703703703
$ws.Range("J9").Value = 0.154320987654321
$ws.Range("O9").Value = 0.02469135802469136
$ws.Range("Q9").Value = 0.2037037037037037
$ws.Range("R9").Value = 0.1358024691358025
$ws.Range("S9").Value = 0.382716049382716
$ws.Range("B10").Value = 0.1124206708975521
$ws.Range("D10").Value = 0.02175883952855848
$ws.Range("F10").Value = 0.06980961015412511
$ws.Range("J10").Value = 0.1278331822302811
$ws.Range("O10").Value = 0.02085222121486854
$ws.Range("Q10").Value = 0.1831368993653672
$ws.Range("R10").Value = 0.0643699002719855
$ws.Range("S10").Value = 0.399818676337262
$ws.Range("G11").Value = 0.1897106109324759
$ws.Range("J11").Value = 0.09646302250803858
$ws.Range("K11").Value = 0.2057877813504823
$ws.Range("L11").Value = 0.4983922829581994
$ws.Range("S11").Value = 0.009646302250803859
$ws.Range("G12").Value = 0.7515151515151515
$ws.Range("J12").Value = 0.1818181818181818
$ws.Range("L12").Value = 0.006060606060606061
$ws.Range("S12").Value = 0.06060606060606061
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3055555555555556
$ws.Range("S13").Value = 0.1111111111111111
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.02645502645502645
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.05291005291005291
$ws.Range("K15").Value = 0.02116402116402116
$ws.Range("M15").Value = 0.01058201058201058
$ws.Range("O15").Value = 0.06878306878306878
$ws.Range("S15").Value = 0.3439153439153439
$ws.Range("F16").Value = 0.04968944099378882
$ws.Range("H16").Value = 0.1677018633540373
$ws.Range("I16").Value = 0.04347826086956522
$ws.Range("J16").Value = 0.4099378881987578
$ws.Range("K16").Value = 0.1490683229813665
$ws.Range("M16").Value = 0.006211180124223602
$ws.Range("O16").Value = 0.04347826086956522
$ws.Range("S16").Value = 0.1304347826086956
$ws.Range("F17").Value = 0.03478260869565217
$ws.Range("H17").Value = 0.1710144927536232
$ws.Range("I17").Value = 0.1043478260869565
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.09565217391304348
$ws.Range("M17").Value = 0.02028985507246377
$ws.Range("N17").Value = 0.002898550724637681
$ws.Range("O17").Value = 0.04927536231884058
$ws.Range("S17").Value = 0.1217391304347826
$ws.Range("H18").Value = 0.1483870967741935
$ws.Range("I18").Value = 0.1225806451612903
$ws.Range("J18").Value = 0.3935483870967742
$ws.Range("K18").Value = 0.1096774193548387
$ws.Range("M18").Value = 0.03225806451612903
$ws.Range("N18").Value = 0.006451612903225806
$ws.Range("O18").Value = 0.05161290322580645
$ws.Range("S18").Value = 0.1354838709677419
$ws.Range("F19").Value = 0.02677824267782427
$ws.Range("H19").Value = 0.197489539748954
$ws.Range("I19").Value = 0.07615062761506276
$ws.Range("J19").Value = 0.3523012552301255
$ws.Range("K19").Value = 0.1263598326359833
$ws.Range("M19").Value = 0.01841004184100418
$ws.Range("O19").Value = 0.06443514644351464
$ws.Range("S19").Value = 0.1380753138075314
